# Insert a new data row at row 583 (shifts existing rows 583-662 down to 584-663)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(583).EntireRow.Insert()

$ws.Range("A583").Value = 6
$ws.Range("B583").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C583").Value = "Metropolitana"
$ws.Range("D583").Value = 44946
$ws.Range("D583").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E583").Value = 13
$ws.Range("F583").Value = 100112044
$ws.Range("G583").Value = "Perejil"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 230
$ws.Range("K583").Value = 14000
$ws.Range("L583").Value = 15000
$ws.Range("M583").Value = 14391
$ws.Range("N583").Value = "$/docena de atados"
$ws.Range("O583").Value = "Región Metropolitana"
$ws.Range("P583").Value = 4797
$ws.Range("Q583").Value = 3
$ws.Range("R583").Value = "Hortaliza"
